$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 'LK644532'
$ws.Range("C5").Value = 'SCL ENTERPRISES LAUNDRY'
$ws.Range("E5").Value = 880
$ws.Range("F5").Value = 'T'
$ws.Range("H5").Value = 45118.04186157407
$ws.Range("I5").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J5").Value = '06/15/23 13:10'
$ws.Range("K5").Value = '06/12/23 19:34'
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = '$880 as of 6/15/2023 11:10:11 AM'
$ws.Range("N5").Value = 880
$ws.Range("O5").Value = 0

# Row 6
$ws.Range("A6").Value = 'L647934'
$ws.Range("C6").Value = 'SB #6'
$ws.Range("E6").Value = 1940
$ws.Range("F6").Value = 'T'
$ws.Range("H6").ClearContents() | Out-Null
$ws.Range("I6").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J6").Value = '04/06/23 22:10'
$ws.Range("K6").Value = '04/06/23 22:05'
$ws.Range("L6").Value = 20
$ws.Range("M6").Value = '$1,940 as of 4/6/2023 8:05:45 PM'
$ws.Range("N6").Value = 1960
$ws.Range("O6").Value = 0

# Row 7
$ws.Range("A7").Value = 'L678988'
$ws.Range("C7").Value = 'PAYELESS MARKET'
$ws.Range("E7").Value = 2440
$ws.Range("F7").Value = 'T'
$ws.Range("H7").Value = 45137.04186157407
$ws.Range("I7").ClearContents() | Out-Null
$ws.Range("J7").Value = '06/15/23 14:31'
$ws.Range("K7").Value = '06/14/23 20:13'
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = '$2,440 as of 6/14/2023 6:13:48 PM'
$ws.Range("N7").Value = 2440
$ws.Range("O7").Value = 0

# Row 8
$ws.Range("A8").Value = 'L688961'
$ws.Range("C8").Value = 'MONA MART'
$ws.Range("E8").Value = 2600
$ws.Range("F8").Value = 'T'
$ws.Range("H8").Value = 45157.04186157407
$ws.Range("I8").ClearContents() | Out-Null
$ws.Range("J8").Value = '06/15/23 13:15'
$ws.Range("K8").Value = '06/14/23 12:41'
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = '$2,600 as of 6/15/2023 11:15:09 AM'
$ws.Range("N8").Value = 2600
$ws.Range("O8").Value = 0

# Row 10
$ws.Range("A10").Value = 'L682801'
$ws.Range("C10").Value = 'SB#5'
$ws.Range("E10").Value = 3440
$ws.Range("F10").Value = 'T'
$ws.Range("H10").Value = 45107.04186157407
$ws.Range("J10").Value = '06/15/23 15:26'
$ws.Range("K10").Value = '06/15/23 15:26'
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = '$3,560 as of 6/14/2023 9:52:10 AM'
$ws.Range("N10").Value = 3540
$ws.Range("O10").Value = 0

# Row 12
$ws.Range("A12").Value = 'L474792'
$ws.Range("C12").Value = 'NICK SHELL SERVICE'
$ws.Range("E12").Value = 4500
$ws.Range("F12").Value = 'T'
$ws.Range("H12").Value = 45111.04186157407
$ws.Range("I12").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J12").Value = '06/14/23 12:44'
$ws.Range("K12").Value = '06/13/23 11:00'
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = '$4,500 as of 6/13/2023 9:00:01 AM'
$ws.Range("N12").Value = 4500
$ws.Range("O12").Value = 0

# Row 13
$ws.Range("A13").Value = 'LK236828'
$ws.Range("C13").Value = 'WORLDWIDE AUTOMOTIVE'
$ws.Range("E13").Value = 5140
$ws.Range("F13").Value = 'T'
$ws.Range("H13").Value = 45109.04186157407
$ws.Range("J13").Value = '06/14/23 20:04'
$ws.Range("K13").Value = '06/14/23 20:04'
$ws.Range("L13").Value = 80
$ws.Range("M13").Value = '$5,140 as of 6/14/2023 6:04:04 PM'
$ws.Range("N13").Value = 5180
$ws.Range("O13").Value = 0

# Row 14
$ws.Range("A14").Value = 'L474817'
$ws.Range("C14").Value = 'SAFETY MARKET'
$ws.Range("E14").Value = 5320
$ws.Range("F14").Value = 'T'
$ws.Range("H14").Value = 45100.04186157407
$ws.Range("J14").Value = '06/15/23 10:14'
$ws.Range("K14").Value = '06/15/23 00:10'
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = '$5,320 as of 6/15/2023 8:14:03 AM'
$ws.Range("N14").Value = 5320
$ws.Range("O14").Value = 0

# Row 15
$ws.Range("A15").Value = 'L662336'
$ws.Range("C15").Value = 'SB#4 MONA MARKET'
$ws.Range("E15").Value = 5700
$ws.Range("F15").Value = 'T'
$ws.Range("H15").Value = 45117.04186157407
$ws.Range("J15").Value = '06/15/23 09:34'
$ws.Range("K15").Value = '06/15/23 09:34'
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = '$5,700 as of 6/15/2023 7:34:59 AM'
$ws.Range("N15").Value = 5780
$ws.Range("O15").Value = 0

# Row 16
$ws.Range("A16").Value = 'LK864765'
$ws.Range("C16").Value = 'SKY LIQUOR'
$ws.Range("E16").Value = 6080
$ws.Range("F16").Value = 'T'
$ws.Range("H16").Value = 45102.04186157407
$ws.Range("I16").ClearContents() | Out-Null
$ws.Range("J16").Value = '06/15/23 15:27'
$ws.Range("K16").Value = '06/15/23 01:51'
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = '$6,080 as of 6/14/2023 11:51:34 PM'
$ws.Range("N16").Value = 6080
$ws.Range("O16").Value = 0

# Row 17
$ws.Range("A17").Value = 'L697589'
$ws.Range("C17").Value = 'S B DISCOUNT MART'
$ws.Range("E17").Value = 6120
$ws.Range("F17").Value = 'T'
$ws.Range("H17").Value = 45097.04186157407
$ws.Range("J17").Value = '06/15/23 12:55'
$ws.Range("K17").Value = '06/15/23 12:55'
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = '$6,120 as of 6/15/2023 10:55:35 AM'
$ws.Range("N17").Value = 6220
$ws.Range("O17").Value = 0

# Row 18
$ws.Range("A18").Value = 'L697590'
$ws.Range("C18").Value = 'S B MARKET ST'
$ws.Range("E18").Value = 6320
$ws.Range("F18").Value = 'T'
$ws.Range("H18").Value = 45108.04186157407
$ws.Range("J18").Value = '06/15/23 15:25'
$ws.Range("K18").Value = '06/15/23 15:25'
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = '$6,440 as of 6/15/2023 8:05:20 AM'
$ws.Range("N18").Value = 6340
$ws.Range("O18").Value = 0

# Row 19
$ws.Range("A19").Value = 'L476340'
$ws.Range("C19").Value = 'DONUT & SANDWICH'
$ws.Range("E19").Value = 6360
$ws.Range("F19").Value = 'T'
$ws.Range("H19").Value = 45121.04186157407
$ws.Range("J19").Value = '06/15/23 15:21'
$ws.Range("K19").Value = '06/15/23 15:21'
$ws.Range("L19").Value = 40
$ws.Range("M19").Value = '$6,540 as of 6/15/2023 9:59:33 AM'
$ws.Range("N19").Value = 6460
$ws.Range("O19").Value = 0

# Row 20
$ws.Range("A20").Value = 'L488595'
$ws.Range("C20").Value = 'N S MART'
$ws.Range("E20").Value = 6360
$ws.Range("F20").Value = 'T'
$ws.Range("H20").Value = 45132.04186157407
$ws.Range("J20").Value = '06/14/23 22:03'
$ws.Range("K20").Value = '06/14/23 22:03'
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = '$6,360 as of 6/14/2023 8:03:30 PM'
$ws.Range("N20").Value = 6420
$ws.Range("O20").Value = 0

# Row 21
$ws.Range("A21").Value = 'L474746'
$ws.Range("C21").Value = 'ZACATES MARKET'
$ws.Range("E21").Value = 6620
$ws.Range("F21").Value = 'T'
$ws.Range("H21").Value = 45121.04186157407
$ws.Range("J21").Value = '06/14/23 18:38'
$ws.Range("K21").Value = '06/14/23 18:38'
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = '$6,620 as of 6/14/2023 4:38:13 PM'
$ws.Range("N21").Value = 6680
$ws.Range("O21").Value = 0

# Row 22
$ws.Range("A22").Value = 'L474761'
$ws.Range("C22").Value = 'BABS MARKET'
$ws.Range("E22").Value = 6920
$ws.Range("F22").Value = 'T'
$ws.Range("H22").Value = 45168.04186157407
$ws.Range("I22").ClearContents() | Out-Null
$ws.Range("J22").Value = '06/15/23 13:03'
$ws.Range("K22").Value = '06/15/23 13:03'
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = '$6,920 as of 6/15/2023 11:03:17 AM'
$ws.Range("N22").Value = 6980
$ws.Range("O22").Value = 0

# Row 23
$ws.Range("A23").Value = 'L688966'
$ws.Range("C23").Value = 'LACON MINI MART'
$ws.Range("E23").Value = 7260
$ws.Range("F23").Value = 'T'
$ws.Range("H23").Value = 45180.04186157407
$ws.Range("J23").Value = '06/15/23 15:31'
$ws.Range("K23").Value = '06/15/23 12:37'
$ws.Range("L23").Value = 20
$ws.Range("M23").Value = '$7,260 as of 6/15/2023 10:37:22 AM'
$ws.Range("N23").Value = 7260
$ws.Range("O23").Value = 0

# Row 24
$ws.Range("A24").Value = 'LK923383'
$ws.Range("C24").Value = 'SAMYS PHONE CARDS'
$ws.Range("E24").Value = 10120
$ws.Range("F24").Value = 'T'
$ws.Range("H24").Value = 45106.04186157407
$ws.Range("I24").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J24").Value = '06/14/23 14:05'
$ws.Range("K24").Value = '06/12/23 22:27'
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = '$10,120 as of 6/12/2023 8:27:39 PM'
$ws.Range("N24").Value = 10220
$ws.Range("O24").Value = 0

# Row 25
$ws.Range("A25").Value = 'L475090'
$ws.Range("C25").Value = 'S.B. 2'
$ws.Range("E25").Value = 12260
$ws.Range("F25").Value = 'T'
$ws.Range("H25").Value = 45114.04186157407
$ws.Range("J25").Value = '06/15/23 15:37'
$ws.Range("K25").Value = '06/15/23 15:37'
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = '$12,300 as of 6/15/2023 9:19:46 AM'
$ws.Range("N25").Value = 12280
$ws.Range("O25").Value = 0

# Row 26
$ws.Range("A26").Value = 'LK891176'
$ws.Range("C26").Value = '98 DISCOUNT STORE'
$ws.Range("E26").Value = 21180
$ws.Range("F26").Value = 'T'
$ws.Range("H26").Value = 45107.04186157407
$ws.Range("J26").Value = '06/15/23 15:30'
$ws.Range("K26").Value = '06/15/23 13:08'
$ws.Range("L26").Value = 60
$ws.Range("M26").Value = '$21,180 as of 6/15/2023 11:08:55 AM'
$ws.Range("N26").Value = 21180
$ws.Range("O26").Value = 0

# Row 27
$ws.Range("A27").Value = 'Total Outstanding Cash Balance:'
$ws.Range("E27").Value = 134140
